$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.193.07'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '1.831.39'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9983'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.38'
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6235'
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9994'
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07377'
$ws.Range('E8').Value = '  -1.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2925'
$ws.Range('E9').Value = '  -0.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.15'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07669'
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('D12').Value = '1.829.21'
$ws.Range('E12').Value = '  -1.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.967'
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6675'
$ws.Range('E14').Value = '  -1.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.59'
$ws.Range('E15').Value = '  -0.68%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008962'
$ws.Range('E16').Value = '  -3.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.872'
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('D18').Value = '29.169.65'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').Value = '2.083.93'
$ws.Range('E19').Value = '  -1.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '235.92'
$ws.Range('E20').Value = '  +1.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.47'
$ws.Range('E21').Value = '  -1.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9991'
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.410'
$ws.Range('E23').Value = '  +3.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9989'
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.23'
$ws.Range('E25').Value = '  -1.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1416'
$ws.Range('E26').Value = '  +1.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.526'
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.66'
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('E29').Value = '  -1.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05797'
$ws.Range('E30').Value = '  +4.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.097'
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.088'
$ws.Range('E32').Value = '  -2.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.206'
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.867'
$ws.Range('E34').Value = '  +0.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7325'
$ws.Range('E35').Value = '  -2.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.140'
$ws.Range('E36').Value = '  -0.59%  '
$ws.Range('E37').Value = '  -2.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.840'
$ws.Range('E38').Value = '  +2.37%  '
$ws.Range('D39').Value = '1.227.75'
$ws.Range('E39').Value = '  -0.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01756'
$ws.Range('E40').Value = '  -2.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.274'
$ws.Range('E41').Value = '  -4.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9174'
$ws.Range('E42').Value = '  +1.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9999'
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.87'
$ws.Range('E44').Value = '  -0.51%  '
$ws.Range('D45').Value = '1.987.79'
$ws.Range('E45').Value = '  -1.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.14'
$ws.Range('E46').Value = '  -2.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5046'
$ws.Range('E47').Value = '  -1.04%  '
$ws.Range('E48').Value = '  -4.17%  '
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4027'
$ws.Range('E49').Value = '  -1.52%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.117'
$ws.Range('E50').Value = '  -0.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1131'
$ws.Range('E51').Value = '  +2.55%  '
